# Unreal Polishing Game Flow
# Rebalance the Attack (column C) stat curve on the Player sheet,
# move the active selection to C12, and set the sheet's page setup
# to A4 portrait (matches the authored workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Player")

$ws.Range("C3").Value  = 40
$ws.Range("C4").Value  = 80
$ws.Range("C5").Value  = 120
$ws.Range("C7").Value  = 200
$ws.Range("C8").Value  = 240
$ws.Range("C9").Value  = 280
$ws.Range("C10").Value = 320
$ws.Range("C11").Value = 360
$ws.Range("C12").Value = 400

[void]$ws.Activate()
[void]$ws.Range("C12").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
